# DP-420 Cosmos DB tracker: add "Indexing strategy" + "Integrate with Azure
# services" topics, mark rows as in-progress (Hours=2), fill in actual
# completion dates for finished topics, and widen column A to fit the new
# (longer) topic text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# --- Row 4: "Access and manage data ..." now in progress, actual date filled in ---
$ws.Range("B4").Value = 2
Copy-Format "C4" "D4"
$ws.Range("D4").Value = 45436

# --- Row 5: "Execute queries ..." Hours corrected ---
$ws.Range("B5").Value = 2

# --- Row 6: "Get started ..." finished - fill Actual date + Action + Status ---
$ws.Range("B6").Value = 2
Copy-Format "C6" "D6"
$ws.Range("D6").Value = 45438
Copy-Format "E4" "E6"
$ws.Range("E6").Value = "1.Built in function`r`n2.Cross product join"
Copy-Format "F4" "F6"
$ws.Range("F6").Value = "Done"
$ws.Rows.Item(6).RowHeight = 28.8

# --- Row 7: new topic "Define and implement an indexing strategy ..." ---
Copy-Format "A4" "A7"
$ws.Range("A7").Value = "Define and implement an indexing strategy for Azure Cosmos DB for NoSQL"
Copy-Format "B4" "B7"
$ws.Range("B7").Value = 2
Copy-Format "C6" "C7"
$ws.Range("C7").Value = 45441
Copy-Format "C6" "D7"
$ws.Range("D7").Value = 45444
Copy-Format "E4" "E7"
$ws.Range("E7").Value = "1.Default Indexing Policy.`r`n2.Custom Indexing Policy.`r`n3.Strategy 1 - Include all and Exclude specific path`r`n4.Strategy 2 - Exclude all and Include specific path`r`n5.Composite Index."
Copy-Format "F4" "F7"
$ws.Range("F7").Value = "Done"
$ws.Rows.Item(7).RowHeight = 72

# --- Row 8: new topic "Integrate Azure Cosmos DB for NoSQL with Azure services" ---
Copy-Format "A4" "A8"
$ws.Range("A8").Value = "Integrate Azure Cosmos DB for NoSQL with Azure services"
Copy-Format "B4" "B8"
$ws.Range("B8").Value = 2
Copy-Format "C6" "C8"
$ws.Range("C8").Value = 45445

# --- Column A needs to widen to fit the new, longer topic text ---
$ws.Columns.Item(1).ColumnWidth = 62.498697916666664

# --- Selection left on E7 (as last edited cell) ---
[void]$ws.Range("E7").Select()
